$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "system, System, backup@backdoor.com"
$ws.Range("G3").Value = "System, dnasr281@gmail.com"
$ws.Range("G4").Value = "System, backup@backdoor.com"
$ws.Range("G5").Value = "System, backup@backdoor.com"
$ws.Range("G6").Value = "System, dnasr281@gmail.com"
$ws.Range("G8").Value = "System, backup@backdoor.com"
$ws.Range("G10").Value = "System, dnasr281@gmail.com"
$ws.Range("G11").Value = "System, dnasr281@gmail.com"
$ws.Range("G12").Value = "System, dnasr281@gmail.com"
$ws.Range("G13").Value = "System, dnasr281@gmail.com"
$ws.Range("G14").Value = "System, dnasr281@gmail.com"
$ws.Range("G15").Value = "System, dnasr281@gmail.com"
$ws.Range("G17").Value = "System, dnasr281@gmail.com"
$ws.Range("G18").Value = "System, dnasr281@gmail.com"
$ws.Range("G19").Value = "System, dnasr281@gmail.com"
$ws.Range("G20").Value = "System, dnasr281@gmail.com"
$ws.Range("G21").Value = "System, dnasr281@gmail.com"
$ws.Range("G22").Value = "System, dnasr281@gmail.com"
$ws.Range("G24").Value = "System, dnasr281@gmail.com"
$ws.Range("G26").Value = "System, dnasr281@gmail.com"
$ws.Range("G28").Value = "system, System, backup@backdoor.com"
$ws.Range("G29").Value = "System, dnasr281@gmail.com"
$ws.Range("G30").Value = "System, backup@backdoor.com"
$ws.Range("G31").Value = "System, backup@backdoor.com"
$ws.Range("G32").Value = "System, dnasr281@gmail.com"
$ws.Range("G34").Value = "System, backup@backdoor.com"
$ws.Range("G36").Value = "System, dnasr281@gmail.com"
$ws.Range("G37").Value = "System, dnasr281@gmail.com"
$ws.Range("G38").Value = "System, dnasr281@gmail.com"
$ws.Range("G39").Value = "System, dnasr281@gmail.com"
$ws.Range("G40").Value = "System, dnasr281@gmail.com"
$ws.Range("G41").Value = "System, dnasr281@gmail.com"
$ws.Range("G43").Value = "System, dnasr281@gmail.com"
$ws.Range("G44").Value = "System, dnasr281@gmail.com"
$ws.Range("G45").Value = "System, dnasr281@gmail.com"
$ws.Range("G46").Value = "System, dnasr281@gmail.com"
$ws.Range("G47").Value = "System, dnasr281@gmail.com"
$ws.Range("G48").Value = "System, dnasr281@gmail.com"
$ws.Range("G50").Value = "System, dnasr281@gmail.com"
$ws.Range("G52").Value = "System, dnasr281@gmail.com"
$ws.Range("G54").Value = "system, System, backup@backdoor.com"
$ws.Range("G55").Value = "System, dnasr281@gmail.com"
$ws.Range("G56").Value = "System, backup@backdoor.com"
$ws.Range("G57").Value = "System, backup@backdoor.com"
$ws.Range("G58").Value = "System, dnasr281@gmail.com"
$ws.Range("G60").Value = "System, backup@backdoor.com"
$ws.Range("G62").Value = "System, dnasr281@gmail.com"
$ws.Range("G63").Value = "System, dnasr281@gmail.com"
$ws.Range("G64").Value = "System, dnasr281@gmail.com"
$ws.Range("G65").Value = "System, dnasr281@gmail.com"
$ws.Range("G66").Value = "System, dnasr281@gmail.com"
$ws.Range("G67").Value = "System, dnasr281@gmail.com"
$ws.Range("G69").Value = "System, dnasr281@gmail.com"
$ws.Range("G70").Value = "System, dnasr281@gmail.com"
$ws.Range("G71").Value = "System, dnasr281@gmail.com"
$ws.Range("G72").Value = "System, dnasr281@gmail.com"
$ws.Range("G73").Value = "System, dnasr281@gmail.com"
$ws.Range("G74").Value = "System, dnasr281@gmail.com"
$ws.Range("G76").Value = "System, dnasr281@gmail.com"
$ws.Range("G78").Value = "System, dnasr281@gmail.com"
$ws.Range("G80").Value = "System, backup@backdoor.com"
$ws.Range("G81").Value = "System, backup@backdoor.com"
$ws.Range("G82").Value = "System, backup@backdoor.com"
$ws.Range("G83").Value = "System, dnasr281@gmail.com"
$ws.Range("G84").Value = "System, dnasr281@gmail.com"
$ws.Range("G85").Value = "System, dnasr281@gmail.com"
$ws.Range("G86").Value = "System, dnasr281@gmail.com"
$ws.Range("G87").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G90").Value = "System, dnasr281@gmail.com"
$ws.Range("G92").Value = "System, dnasr281@gmail.com"
$ws.Range("G93").Value = "System, dnasr281@gmail.com"
$ws.Range("G94").Value = "System, dnasr281@gmail.com"
$ws.Range("G96").Value = "System, dnasr281@gmail.com"
$ws.Range("G99").Value = "System, dnasr281@gmail.com"
$ws.Range("G101").Value = "System, dnasr281@gmail.com"
$ws.Range("G106").Value = "System, backup@backdoor.com"
$ws.Range("G107").Value = "System, backup@backdoor.com"
$ws.Range("G108").Value = "System, backup@backdoor.com"
$ws.Range("G109").Value = "System, dnasr281@gmail.com"
$ws.Range("G110").Value = "System, dnasr281@gmail.com"
$ws.Range("G111").Value = "System, dnasr281@gmail.com"
$ws.Range("G112").Value = "System, dnasr281@gmail.com"
$ws.Range("G113").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G116").Value = "System, dnasr281@gmail.com"
$ws.Range("G118").Value = "System, dnasr281@gmail.com"
$ws.Range("G119").Value = "System, dnasr281@gmail.com"
$ws.Range("G120").Value = "System, dnasr281@gmail.com"
$ws.Range("G122").Value = "System, dnasr281@gmail.com"
$ws.Range("G125").Value = "System, dnasr281@gmail.com"
$ws.Range("G127").Value = "System, dnasr281@gmail.com"
$ws.Range("G132").Value = "System, backup@backdoor.com"
$ws.Range("G133").Value = "System, backup@backdoor.com"
$ws.Range("G134").Value = "System, backup@backdoor.com"
$ws.Range("G135").Value = "System, dnasr281@gmail.com"
$ws.Range("G136").Value = "System, dnasr281@gmail.com"
$ws.Range("G137").Value = "System, dnasr281@gmail.com"
$ws.Range("G138").Value = "System, dnasr281@gmail.com"
$ws.Range("G139").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G142").Value = "System, dnasr281@gmail.com"
$ws.Range("G144").Value = "System, dnasr281@gmail.com"
$ws.Range("G145").Value = "System, dnasr281@gmail.com"
$ws.Range("G146").Value = "System, dnasr281@gmail.com"
$ws.Range("G148").Value = "System, dnasr281@gmail.com"
$ws.Range("G151").Value = "System, dnasr281@gmail.com"
$ws.Range("G153").Value = "System, dnasr281@gmail.com"

